$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 179, pushing existing rows 179-192 down to 183-196
$ws.Rows("179:182").Insert()

# Row 179
$ws.Range("A179").Value = 1
$ws.Range("B179").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C179").Value = 'Arica y Parinacota'
$ws.Range("D179").Value = 45013
$ws.Range("E179").Value = 15
$ws.Range("F179").Value = 'Fruta'
$ws.Range("G179").Value = 100104
$ws.Range("H179").Value = 'Frutos de pepita'
$ws.Range("I179").Value = 100104002
$ws.Range("J179").Value = 'Manzana'
$ws.Range("K179").Value = 'Ambrosia'
$ws.Range("L179").Value = 'Segunda'
$ws.Range("M179").Value = 260
$ws.Range("N179").Value = 20000
$ws.Range("O179").Value = 22000
$ws.Range("P179").Value = 20923
$ws.Range("Q179").Value = '$/caja 18 kilos empedrada'
$ws.Range("R179").Value = 'Provincia de Curicó'
$ws.Range("S179").Value = 1162
$ws.Range("T179").Value = 18

# Row 180
$ws.Range("A180").Value = 1
$ws.Range("B180").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C180").Value = 'Arica y Parinacota'
$ws.Range("D180").Value = 45013
$ws.Range("E180").Value = 15
$ws.Range("F180").Value = 'Fruta'
$ws.Range("G180").Value = 100104
$ws.Range("H180").Value = 'Frutos de pepita'
$ws.Range("I180").Value = 100104002
$ws.Range("J180").Value = 'Manzana'
$ws.Range("K180").Value = 'Fuji royal'
$ws.Range("L180").Value = 'Segunda'
$ws.Range("M180").Value = 250
$ws.Range("N180").Value = 20000
$ws.Range("O180").Value = 22000
$ws.Range("P180").Value = 21040
$ws.Range("Q180").Value = '$/caja 18 kilos empedrada'
$ws.Range("R180").Value = 'Provincia de Curicó'
$ws.Range("S180").Value = 1169
$ws.Range("T180").Value = 18

# Row 181
$ws.Range("A181").Value = 1
$ws.Range("B181").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C181").Value = 'Arica y Parinacota'
$ws.Range("D181").Value = 45013
$ws.Range("E181").Value = 15
$ws.Range("F181").Value = 'Fruta'
$ws.Range("G181").Value = 100104
$ws.Range("H181").Value = 'Frutos de pepita'
$ws.Range("I181").Value = 100104002
$ws.Range("J181").Value = 'Manzana'
$ws.Range("K181").Value = 'Granny Smith'
$ws.Range("L181").Value = 'Segunda'
$ws.Range("M181").Value = 220
$ws.Range("N181").Value = 20000
$ws.Range("O181").Value = 22000
$ws.Range("P181").Value = 20909
$ws.Range("Q181").Value = '$/caja 18 kilos empedrada'
$ws.Range("R181").Value = 'Provincia de Curicó'
$ws.Range("S181").Value = 1162
$ws.Range("T181").Value = 18

# Row 182
$ws.Range("A182").Value = 1
$ws.Range("B182").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C182").Value = 'Arica y Parinacota'
$ws.Range("D182").Value = 45013
$ws.Range("E182").Value = 15
$ws.Range("F182").Value = 'Fruta'
$ws.Range("G182").Value = 100104
$ws.Range("H182").Value = 'Frutos de pepita'
$ws.Range("I182").Value = 100104002
$ws.Range("J182").Value = 'Manzana'
$ws.Range("K182").Value = 'Royal Gala'
$ws.Range("L182").Value = 'Segunda'
$ws.Range("M182").Value = 250
$ws.Range("N182").Value = 20000
$ws.Range("O182").Value = 22000
$ws.Range("P182").Value = 20800
$ws.Range("Q182").Value = '$/caja 18 kilos empedrada'
$ws.Range("R182").Value = 'Provincia de Curicó'
$ws.Range("S182").Value = 1156
$ws.Range("T182").Value = 18
